# Apply updated crypto market data (prices and 1h volume change) per row.
# Numeric-looking "Price" values are written with a leading apostrophe so Excel
# keeps them as plain text (matching the source data, which stores everything as text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.568.13'
$ws.Range("E2").Value = '  +5.79%  '

$ws.Range("D3").Value = '2.398.21'
$ws.Range("E3").Value = '  +3.97%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = "'113.92"
$ws.Range("E5").Value = '  +8.06%  '

$ws.Range("D6").Value = "'319.85"
$ws.Range("E6").Value = '  +3.00%  '

$ws.Range("E7").Value = '  +1.13%  '

$ws.Range("E8").Value = '  -0.31%  '

$ws.Range("E9").Value = '  +3.29%  '

$ws.Range("D10").Value = "'42.29"
$ws.Range("E10").Value = '  +6.35%  '

$ws.Range("D11").Value = "'0.0930"
$ws.Range("E11").Value = '  +2.44%  '

$ws.Range("D12").Value = "'8.70"
$ws.Range("E12").Value = '  +4.93%  '

$ws.Range("E13").Value = '  +2.59%  '

$ws.Range("E14").Value = '  +1.15%  '

$ws.Range("D15").Value = "'15.89"
$ws.Range("E15").Value = '  +3.60%  '

$ws.Range("D16").Value = '2.765.42'
$ws.Range("E16").Value = '  +5.28%  '

$ws.Range("D17").Value = '2.382.36'
$ws.Range("E17").Value = '  +3.40%  '

$ws.Range("D18").Value = '45.546.74'
$ws.Range("E18").Value = '  +6.13%  '

$ws.Range("D19").Value = "'7.51"
$ws.Range("E19").Value = '  +2.44%  '

$ws.Range("E20").Value = '  +3.13%  '

$ws.Range("D21").Value = "'13.34"
$ws.Range("E21").Value = '  -1.69%  '

$ws.Range("D22").Value = "'74.76"
$ws.Range("E22").Value = '  +1.71%  '

$ws.Range("D23").Value = "'3.58"
$ws.Range("E23").Value = '  +3.98%  '

$ws.Range("D24").Value = "'265.10"
$ws.Range("E24").Value = '  -0.76%  '

$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = '  +4.29%  '

$ws.Range("E26").Value = '  -1.02%  '

$ws.Range("D27").Value = "'7.68"
$ws.Range("E27").Value = '  +5.33%  '

$ws.Range("D28").Value = "'11.35"
$ws.Range("E28").Value = '  +3.23%  '

$ws.Range("E29").Value = '  +2.90%  '

$ws.Range("D30").Value = "'39.53"
$ws.Range("E30").Value = '  +4.85%  '

$ws.Range("D31").Value = "'22.74"
$ws.Range("E31").Value = '  +2.03%  '

$ws.Range("D32").Value = "'0.0972"
$ws.Range("E32").Value = '  +12.88%  '

$ws.Range("D33").Value = "'173.05"
$ws.Range("E33").Value = '  +4.78%  '

$ws.Range("E34").Value = '  +3.70%  '

$ws.Range("E35").Value = '  +1.45%  '

$ws.Range("E36").Value = '  +7.59%  '

$ws.Range("E37").Value = '  +6.23%  '

$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = "'3.09"
$ws.Range("E38").Value = '  +8.41%  '

$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").Value = "'4.11"
$ws.Range("E39").Value = '  +13.03%  '

$ws.Range("E40").Value = '  +4.03%  '

$ws.Range("E41").Value = '  +12.19%  '

$ws.Range("B42").Value = 'BitcoinSV'
$ws.Range("C42").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D42").Value = "'103.32"
$ws.Range("E42").Value = '  -4.16%  '

$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").Value = "'13.82"
$ws.Range("E43").Value = '  +12.33%  '

$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = "'0.241"
$ws.Range("E44").Value = '  +5.20%  '

$ws.Range("D45").Value = "'71.55"
$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("D46").Value = "'88.23"
$ws.Range("E46").Value = '  +14.78%  '

$ws.Range("E47").Value = '  -0.46%  '

$ws.Range("D48").Value = "'5.81"
$ws.Range("E48").Value = '  +12.32%  '

$ws.Range("D49").Value = "'116.30"
$ws.Range("E49").Value = '  +4.07%  '

$ws.Range("D50").Value = "'9.50"
$ws.Range("E50").Value = '  +8.58%  '

$ws.Range("D51").Value = '1.678.26'
$ws.Range("E51").Value = '  -2.60%  '
